# Fixing the big mistake: correct the Total (B) and Community (D) values
# for rows 2-13 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 2;  B = 14444.12194975002;  D = 666.4819821333334  }
    @{ Row = 3;  B = 13530.57041270002;  D = 629.4793040000001  }
    @{ Row = 4;  B = 14472.11968053336;  D = 685.04545515       }
    @{ Row = 5;  B = 13996.37456801669;  D = 637.5027093833334  }
    @{ Row = 6;  B = 14531.34629450002;  D = 690.6329750666666  }
    @{ Row = 7;  B = 13997.48681761669;  D = 673.12289          }
    @{ Row = 8;  B = 14430.94576628336;  D = 647.9368861833333  }
    @{ Row = 9;  B = 14499.41614671669;  D = 694.7404019333333  }
    @{ Row = 10; B = 14057.64996260002;  D = 646.61009635       }
    @{ Row = 11; B = 14470.00854215002;  D = 674.1076421333332  }
    @{ Row = 12; B = 14020.68100321669;  D = 629.6293876166667  }
    @{ Row = 13; B = 14019.26579695003;  D = 653.7429827833333  }
)

foreach ($entry in $values) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
    $ws.Cells.Item($entry.Row, 4).Value = $entry.D
}
